$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. "Datos faltantes:" -> underline "Datos faltantes", leave ":" plain
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Datos faltantes:", $false, $false, $false, $false, $false, `
                 $true, 1, $false, "", 0) | Out-Null
if ($r.Find.Found) {
    $wordEnd = $r.End - 1
    $textRange = $d.Range($r.Start, $wordEnd)
    $textRange.Font.Underline = 1
}

# ---------------------------------------------------------------------
# 2. "Datos descompensados:" -> strikethrough whole paragraph (mark + run)
# ---------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("Datos descompensados:", $false, $false, $false, $false, `
                  $false, $true, 1, $false, "", 0) | Out-Null
if ($r2.Find.Found) {
    $para = $r2.Paragraphs(1)
    $para.Range.Font.StrikeThrough = 1
}

# ---------------------------------------------------------------------
# 3. "GridSearch" -> strikethrough just the word
# ---------------------------------------------------------------------
$r3 = $d.Content
$r3.Find.Execute("GridSearch", $false, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0) | Out-Null
if ($r3.Find.Found) {
    $r3.Font.StrikeThrough = 1
}

# ---------------------------------------------------------------------
# 4. Remove the _GoBack bookmark
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 5. Merge the split "https://tow" + "a" + "rdsdatascience.com/..." runs
#    of the MLBOX hyperlink into a single run with the full URL text
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Hyperlinks.Count; $i++) {
    $h = $d.Hyperlinks($i)
    if ($h.TextToDisplay -eq "https://towardsdatascience.com/automl-in-python-an-overview-of-the-mlbox-package-208118a7fe5") {
        $h.TextToDisplay = "https://towardsdatascience.com/automl-in-python-an-overview-of-the-mlbox-package-208118a7fe5"
        break
    }
}

Write-Output "done"
